# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 12 (pushing existing rows 12-39 down to 13-40),
# matching the new Magnum / Región Metropolitana observation added to the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 12; Excel carries the D-column date
# style (s="2") down from the row above automatically.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new record's values.
$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C12").Value = "Ñuble"
$ws.Range("D12").Value = 44526
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 100112031
$ws.Range("G12").Value = "Poroto verde"
$ws.Range("H12").Value = "Magnum"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 29000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 29500
$ws.Range("N12").Value = "$/saco 25 kilos"
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 1180
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
